# Weekly update: a new price record is inserted at the top of the
# Espinaca / Vega Modelo de Temuco series (row 94), every existing
# record from row 94..216 shifts down by one row, and the record that
# used to be the last one (row 216) is preserved by appending a new
# row 217.
#
# Columns that move together as one "record" (everything that is not
# constant across the whole sheet): D (Fecha), I (Calidad), J (Volumen),
# K (Precio minimo), L (Precio maximo), M (Precio promedio ponderado),
# N (Unidad de comercializacion), O (Origen), P (Precio $/Kg).
# Columns A,B,C,E,F,G,H,Q,R are identical on every data row, so row 217
# can just copy them straight from row 216.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 94
$lastRow = 216
$newLastRow = 217

# Column numbers for the "moving" fields: D,I,J,K,L,M,N,O,P
$movingCols = @(4, 9, 10, 11, 12, 13, 14, 15, 16)

# 1) Snapshot every current value we will need, BEFORE any write,
#    since the shift reads row r-1's *original* value while writing row r.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @{}
    foreach ($c in $movingCols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# Also snapshot the static columns of the last row, to build the new
# appended row 217.
$staticCols = @(1, 2, 3, 5, 6, 7, 8, 17, 18)
$staticVals = @{}
foreach ($c in $staticCols) {
    $staticVals[$c] = $ws.Cells.Item($lastRow, $c).Value2
}

# The Fecha (date) column carries a custom date number-format; remember
# it so the freshly appended row 217 can reuse it.
$dateFormat = $ws.Cells.Item($lastRow, 4).NumberFormat

# 2) Shift every record down by one row: new row r (for r = 95 .. 217)
#    gets the values that used to live in row r-1 (94 .. 216). We only
#    ever read from the $snapshot hashtable (captured in step 1), never
#    from the live sheet, so the write order does not matter.
for ($dstRow = $newLastRow; $dstRow -ge ($firstRow + 1); $dstRow--) {
    $src = $snapshot[$dstRow - 1]
    foreach ($c in $movingCols) {
        $ws.Cells.Item($dstRow, $c).Value = $src[$c]
    }
}

# 3) New row 217 gets the static columns copied from old row 216 (the
#    moving columns for row 217 were already written above, since
#    dstRow = 217 pulls from snapshot[216]).
foreach ($c in $staticCols) {
    $ws.Cells.Item($newLastRow, $c).Value = $staticVals[$c]
}
$ws.Cells.Item($newLastRow, 4).NumberFormat = $dateFormat

# 4) Row 94 becomes the brand-new record. Its static columns already
#    hold the right (unchanged) values, and I/N/O keep their previous
#    value (Primera / $/docena de atados / Region de La Araucania), so
#    only D, J, K, L, M, P need the new figures.
$ws.Cells.Item($firstRow, 4).Value = 44902
$ws.Cells.Item($firstRow, 10).Value = 55
$ws.Cells.Item($firstRow, 11).Value = 10000
$ws.Cells.Item($firstRow, 12).Value = 10000
$ws.Cells.Item($firstRow, 13).Value = 10000
$ws.Cells.Item($firstRow, 16).Value = 3333

# 5) Make sure the used-range dimension now covers the new row.
$ws.Range("A1:R217").Select() | Out-Null
